# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: subscript-three character used in PEPE price (U+2083)
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.882.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.919.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "639.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.829.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.378.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.913"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +6.92%  "
$ws.Range("E27").Value = "  +4.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "611.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.058.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  +6.52%  "
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0$sub3" + "0705"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0423"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("E48").Value = "  +12.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.29%  "
